# dmb ran blue tank titrations 0505
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Find the next empty row after the existing data block (row 32 -> 33).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Append the new titration data row, carrying the column D "% off" formula
# down the same way it was filled for the previous rows.
$ws.Cells.Item($newRow, 1).Value = 20210505
$ws.Cells.Item($newRow, 2).Value = 2229.0880000000002
$ws.Cells.Item($newRow, 3).Value = 2224.4699999999998
$ws.Cells.Item($newRow, 4).Formula = "=100*(B" + $newRow + "-C" + $newRow + ")/C" + $newRow
$ws.Cells.Item($newRow, 5).Value = 180
$ws.Cells.Item($newRow, 6).Value = "CRM opened 20210418"

$ws.Cells.Item($newRow, 4).Select()
